$d = $word.ActiveDocument
$tab = [char]9
$vtab = [char]11

# --- Contact info block ---
$d.Content.Find.Execute("New York", $true, $false, $false, $false, $false, $true, 1, $false, "New York, USA", 2)
$d.Content.Find.Execute("alex@gmail.com", $true, $false, $false, $false, $false, $true, 1, $false, "alex@example.com", 2)
$d.Content.Find.Execute("32929823", $true, $false, $false, $false, $false, $true, 1, $false, "00000000", 2)

# --- Experiences section ---
# spaceX, <tab>  ->  Google,   (tab removed, text changed)
$d.Content.Find.Execute("spaceX, $tab", $true, $false, $false, $false, $false, $true, 1, $false, "Google, ", 2)

# complete the first experience's Lorem Ipsum paragraph (scoped to that paragraph only,
# since "readable Engl" is also a substring further up in "Career Profile")
$p9 = $d.Paragraphs.Item(9).Range
$p9.Find.Execute("readable Engl", $true, $false, $false, $false, $false, $true, 1, $false, "readable English.", 2)

# Google, <tab>  ->  Tesla,   (tab removed, text changed) -- only remaining match now
$d.Content.Find.Execute("Google, $tab", $true, $false, $false, $false, $false, $true, 1, $false, "Tesla, ", 2)

# dates for the (now) Tesla entry
$d.Content.Find.Execute("2017-2019", $true, $false, $false, $false, $false, $true, 1, $false, "2016-2018", 2)

# truncate the second experience's Lorem Ipsum paragraph (scoped to that paragraph)
$p10 = $d.Paragraphs.Item(10).Range
$oldText = "It is a long established fact that a reader will be distracted by the readable content of a page when looking at its layout. The point of using Lorem Ipsum is that it has a more-or-less normal distribution of letters, as opposed to using '"
$newText = "It is a long established fact that a reader will be distracted by the readable content of a page when looking at its layout. The point of using Lorem Ipsum is that it has a more-or-less."
$p10.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# --- New Education section ---
$last = $d.Paragraphs.Last

# 1) blank separator paragraph
$last.Range.InsertParagraphBefore()

# 2) "Education" heading
$last.Range.InsertParagraphBefore()
$headingPara = $last.Previous
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "Education"

# 3) "Example University" / 2015-2017
$last.Range.InsertParagraphBefore()
$uniPara = $last.Previous
$uniRange = $uniPara.Range
$uniStart = $uniRange.Start
$uniRange.Text = "Example University" + $vtab + "2015-2017"
$uniTitleLen = ("Example University").Length
$uniBoldRange = $d.Range($uniStart, $uniStart + $uniTitleLen)
$uniBoldRange.Bold = 1

# 4) "Example School" / 2012-2014
$last.Range.InsertParagraphBefore()
$schoolPara = $last.Previous
$schoolRange = $schoolPara.Range
$schoolStart = $schoolRange.Start
$schoolRange.Text = "Example School" + $vtab + "2012-2014"
$schoolTitleLen = ("Example School").Length
$schoolBoldRange = $d.Range($schoolStart, $schoolStart + $schoolTitleLen)
$schoolBoldRange.Bold = 1

Write-Host "done"
